$wb = $excel.ActiveWorkbook

# Rename the "DiagnosisDuplicates" sheet to "DiagnosisRedundancies"
$ws = $wb.Worksheets.Item("DiagnosisDuplicates")
$ws.Name = "DiagnosisRedundancies"

# Update the column header text on that sheet from "IsLikelyDuplicate" to "IsLikelyRedundant"
$ws.Range("B4:B11").Value = "IsLikelyRedundant"
